# Append a new user-detail record (Ewan Marsh) as row 33 on Sheet1,
# mirroring the other existing rows' layout/values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = 110032
$ws.Range("B33").Value = 9317596770
$ws.Range("C33").Value = "Ewan Marsh"
$ws.Range("D33").Value = "ewan.marsh@xyz.com"
$ws.Range("E33").Value = 818876433
$ws.Range("F33").Value = "ACT"
$ws.Range("G33").Value = "eng"
$ws.Range("H33").Value = "PWD"
$ws.Range("I33").Value = $true
$ws.Range("J33").Value = "superadmin"
$ws.Range("K33").Value = "now()"
$ws.Range("L33").Value = "now()"

# Reset the view: scroll back to the top-left and select the unused
# columns to the right of the data (M:XFD), as was left selected after
# the edit.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Columns("M:XFD").Select()
